$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at the top for the "Common Technologies" header
$ws.Rows.Item(1).Insert()

# 2. Set the new header cell value and give it the yellow-fill / non-bold style
$ws.Cells.Item(1,1).Value = "Common Technologies"
$ws.Cells.Item(1,1).Interior.Color = 65535
$ws.Cells.Item(1,1).Font.Bold = $false

# 3. Fill in the previously-empty "P" marker cells in column B (rows shifted by the insert)
$ws.Cells.Item(2,2).Value = "P"
$ws.Cells.Item(3,2).Value = "P"
$ws.Cells.Item(4,2).Value = "P"
$ws.Cells.Item(5,2).Value = "Didn’t need to"
$ws.Cells.Item(6,2).Value = "P"
$ws.Cells.Item(7,2).Value = "P"
$ws.Cells.Item(8,2).Value = "P"
$ws.Cells.Item(9,2).Value = "P"
$ws.Cells.Item(10,2).Value = "P"
$ws.Cells.Item(11,2).Value = "P"
$ws.Cells.Item(12,2).Value = "P"

# Ensure column B entries use the existing "P"-style formatting (copy format from B3 template)
$ws.Cells.Item(3,2).Copy()
$destB = $ws.Range("B2,B3,B4,B5,B6,B7,B8,B9,B10,B11,B12")
$destB.PasteSpecial(-4122)

# 4. Append the new "Localization" section at the bottom of the sheet
$ws.Cells.Item(38,1).Value = "Localization"
$ws.Cells.Item(38,1).Interior.Color = 65535
$ws.Cells.Item(38,1).Font.Bold = $false

$ws.Cells.Item(39,1).Value = "stats_l_english.yml"
$ws.Cells.Item(40,1).Value = "technology_sharing_l_english.yml"
$ws.Cells.Item(41,1).Value = "traits_l_english.yml"
$ws.Cells.Item(42,1).Value = "unit_l_english.yml"
$ws.Cells.Item(43,1).Value = "victory_points_l_english.yml"
$ws.Cells.Item(44,1).Value = "war_l_english.yml"
$ws.Cells.Item(45,1).Value = "bookmarks_l_english.yml"
$ws.Cells.Item(46,1).Value = "countries_cosmetic_l_english.yml"
$ws.Cells.Item(47,1).Value = "countries_l_english.yml"
$ws.Cells.Item(48,1).Value = "equip_air_l_english.yml"
$ws.Cells.Item(49,1).Value = "equip_naval_l_english.yml"
$ws.Cells.Item(50,1).Value = "equipment_l_english.yml"
$ws.Cells.Item(51,1).Value = "events_l_english.yml"
$ws.Cells.Item(52,1).Value = "frontend_l_english.yml"
$ws.Cells.Item(53,1).Value = "ideas_l_english.yml"
$ws.Cells.Item(54,1).Value = "loading_tips_l_english.yml"
$ws.Cells.Item(55,1).Value = "modifiers_l_english.yml"
$ws.Cells.Item(56,1).Value = "nef_council_l_english.yml"
$ws.Cells.Item(57,1).Value = "nef_despdefense_l_english.yml"
$ws.Cells.Item(58,1).Value = "nef_eventcountry_l_english.yml"
$ws.Cells.Item(59,1).Value = "nef_eventelection_l_english.yml"
$ws.Cells.Item(60,1).Value = "nef_eventnews_l_english.yml"
$ws.Cells.Item(61,1).Value = "nef_factions_l_english.yml"
$ws.Cells.Item(62,1).Value = "nef_focus_resistance_l_english.yml"
$ws.Cells.Item(63,1).Value = "nef_focus_totalitarian_l_english.yml"
$ws.Cells.Item(64,1).Value = "nef_ideas_l_english.yml"
$ws.Cells.Item(65,1).Value = "nef_other_l_english.yml"
$ws.Cells.Item(66,1).Value = "parties_l_english.yml"
$ws.Cells.Item(67,1).Value = "research_l_english.yml"
$ws.Cells.Item(68,1).Value = "state_names_l_english.yml"

# Apply the grey body style (same as used for rows 14-37) to the new Localization rows
$ws.Cells.Item(14,1).Copy()
$ws.Range("A39:A68").PasteSpecial(-4122)

# 5. Update the active selection to match the final workbook state
$ws.Range("E11").Select()

Write-Host "done"
